# Crystal RC values updated. New pcb images added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C10, C11 capacitor value: 25pF -> 20pF
$ws.Range("B29").Value = "20pF"

# Update R9 resistor value: 400R -> 330k
$ws.Range("B34").Value = "330k"

# Update tactile switch row: Parts list now shows actual designators,
# and the Description now calls out those designators too.
$ws.Range("D53").Value = "S2, S3, S4, S6"
$ws.Range("E53").Value = "5x5mm push button SMD (USER, NRST, PWRKEY, RESET)"

# Reflect the user's final selection/scroll position on the sheet
$ws.Activate()
$ws.Range("B29").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
